$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9328075051307678
$ws.Range("B1").Value = 1.510351061820984
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.360382199287415
$ws.Range("E1").Value = 1.34572970867157
